$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 data (7th iteration)
$ws.Range("A8").Value = 7
$ws.Range("C8").Value = "Joel Rosario, Eric Cancel, Manuel Franco, Andre Shivnarie Worrie, Reylu Gutierres"
$ws.Range("D8").Value = "track_id, race_number, race_date, jockey, program_number"
$ws.Range("E8").Value = "race_date, latitude, longitude, trakus_index, program_number, distance_id, run_up_distance, purse, post_time, odds"
$ws.Range("F8").Value = "Lasso and Ridge Regression performed for comparison. Ultimately decided that unscaled Linear Regression yielded the best results. "

# Match formatting of the preceding note rows (wrap text + row height)
$ws.Range("C8:F8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 85

# Update the active selection to reflect where editing left off
$ws.Range("D11").Select()
